$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# Columns A-D hold text values in this sheet (dates/times/weekday/week are
# stored as literal strings, not parsed as Excel dates/numbers). "16:05:26"
# and "Tuesday" round-trip as text naturally, but "2023-06-06" and "23" look
# like a date/number to Excel's auto-detection, so force text entry for
# those via a temporary "@" number format, then clear the format residue so
# the cell ends up with no explicit style - matching the rest of the sheet.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-06"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "16:05:26"
$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "23"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 120497
$ws.Cells.Item($row, 6).Value = 134182
$ws.Cells.Item($row, 7).Value = 159648
$ws.Cells.Item($row, 8).Value = 130517
$ws.Cells.Item($row, 9).Value = 175163
$ws.Cells.Item($row, 10).Value = 112452
$ws.Cells.Item($row, 11).Value = 200250
$ws.Cells.Item($row, 12).Value = 219855
$ws.Cells.Item($row, 13).Value = 172404
$ws.Cells.Item($row, 14).Value = 119505
$ws.Cells.Item($row, 15).Value = 38364
$ws.Cells.Item($row, 16).Value = 34627
$ws.Cells.Item($row, 17).Value = 50484
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36590
$ws.Cells.Item($row, 20).Value = -1

$wb.Save()
